$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (in-place rich text run edits) ---
# A8: "Volume 32   Number  9" -> "...10"  (replace the trailing run "9")
$ws.Range("A8").Characters(21, 1).Text = "10"
# C9: "Report Covering the Week  2/24/2025  Through  3/2/2025"
#     -> "...3/3/2025  Through  3/9/2025"
# Replace the later run first so the earlier offset is not shifted by a length change.
$ws.Range("C9").Characters(47, 8).Text = "3/9/2025"
$ws.Range("C9").Characters(27, 9).Text = "3/3/2025"

# --- C15: numeric 2 -> text "0" (shared string), keep style (copy format from D15) ---
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null

# --- C27: numeric 2 -> text "0" (shared string), keep style (copy format from D27) ---
$ws.Range("C27").Value = "'0"
$ws.Range("D27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null

# --- D29/E29, D30/E30: text "N/A"/"***.*"-> numeric; copy numeric formats from D20/E20 ---
$ws.Range("D20").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Value = 1
$ws.Range("E20").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = -100

$ws.Range("D20").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("D30").Value = 1
$ws.Range("E20").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = -100

# --- Remaining simple numeric value updates ---
# Row 15
$ws.Range("F15").Value = 8
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 9
$ws.Range("K15").Value = 80
$ws.Range("L15").Value = 350
$ws.Range("M15").Value = 800
$ws.Range("N15").Value = 80

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 15
$ws.Range("I16").Value = 57
$ws.Range("J16").Value = 58
$ws.Range("K16").Value = -1.724137931034
$ws.Range("L16").Value = -42.424242424242
$ws.Range("M16").Value = 67.647058823529
$ws.Range("N16").Value = -87.276785714285

# Row 17
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = -55.555555555555
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = -48.333333333333
$ws.Range("I17").Value = 89
$ws.Range("J17").Value = 110
$ws.Range("K17").Value = -19.090909090909
$ws.Range("L17").Value = -23.275862068965
$ws.Range("M17").Value = 111.904761904762
$ws.Range("N17").Value = -27.642276422764

# Row 18
$ws.Range("D18").Value = 12
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 37
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = 19.354838709677
$ws.Range("I18").Value = 78
$ws.Range("J18").Value = 77
$ws.Range("K18").Value = 1.298701298701
$ws.Range("L18").Value = -17.021276595744
$ws.Range("M18").Value = 6.849315068493
$ws.Range("N18").Value = -85.474860335195

# Row 19
$ws.Range("D19").Value = 41
$ws.Range("E19").Value = -12.195121951219
$ws.Range("F19").Value = 136
$ws.Range("G19").Value = 133
$ws.Range("H19").Value = 2.255639097744
$ws.Range("I19").Value = 349
$ws.Range("J19").Value = 454
$ws.Range("K19").Value = -23.127753303964
$ws.Range("L19").Value = -20.681818181818
$ws.Range("M19").Value = -15.291262135922
$ws.Range("N19").Value = -81.063483450895

# Row 20
$ws.Range("J20").Value = 9
$ws.Range("K20").Value = -55.555555555555
$ws.Range("M20").Value = -33.333333333333
$ws.Range("N20").Value = -94.736842105263

# Row 21
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 76
$ws.Range("E21").Value = -26.315789473684
$ws.Range("F21").Value = 237
$ws.Range("G21").Value = 251
$ws.Range("H21").Value = -5.577689243027
$ws.Range("I21").Value = 586
$ws.Range("J21").Value = 715
$ws.Range("K21").Value = -18.041958041958
$ws.Range("L21").Value = -22.996057818659
$ws.Range("M21").Value = 3.169014084507
$ws.Range("N21").Value = -80.691927512355

# Row 22
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 33.333333333333
$ws.Range("F22").Value = 12
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = 9.090909090909
$ws.Range("I22").Value = 28
$ws.Range("J22").Value = 36
$ws.Range("K22").Value = -22.222222222222
$ws.Range("M22").Value = 0

# Row 24
$ws.Range("C24").Value = 73
$ws.Range("D24").Value = 84
$ws.Range("E24").Value = -13.095238095238
$ws.Range("F24").Value = 302
$ws.Range("G24").Value = 346
$ws.Range("H24").Value = -12.716763005780
$ws.Range("I24").Value = 731
$ws.Range("J24").Value = 784
$ws.Range("K24").Value = -6.760204081632
$ws.Range("L24").Value = 11.263318112633
$ws.Range("M24").Value = -10.196560196560

# Row 25
$ws.Range("C25").Value = 63
$ws.Range("D25").Value = 76
$ws.Range("E25").Value = -17.105263157894
$ws.Range("F25").Value = 256
$ws.Range("G25").Value = 300
$ws.Range("H25").Value = -14.666666666666
$ws.Range("I25").Value = 641
$ws.Range("J25").Value = 662
$ws.Range("K25").Value = -3.172205438066
$ws.Range("L25").Value = 3.387096774193

# Row 26
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = -16
$ws.Range("F26").Value = 69
$ws.Range("G26").Value = 82
$ws.Range("H26").Value = -15.853658536585
$ws.Range("I26").Value = 192
$ws.Range("J26").Value = 193
$ws.Range("K26").Value = -0.518134715025
$ws.Range("L26").Value = 18.518518518518
$ws.Range("M26").Value = 65.517241379310

# Row 27
$ws.Range("I27").Value = 11
$ws.Range("K27").Value = 57.142857142857
$ws.Range("L27").Value = 120

# Row 28
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 7
$ws.Range("F28").Value = 20
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 40
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = 14.285714285714
$ws.Range("L28").Value = -2.439024390243

# Row 29
$ws.Range("J29").Value = 3

# Row 30
$ws.Range("J30").Value = 3

# Row 31
$ws.Range("D31").Value = 3
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = -75
$ws.Range("J31").Value = 4
$ws.Range("K31").Value = -25
$ws.Range("L31").Value = -25

